$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Data" to "Summary"
$ws.Name = "Summary"

# Shift the header row (currently row 5: Micro/SMEs/MSMEs) and the data
# label row (currently row 6: "Enterprises (% of total)") down so a new
# "Source Type" line and the new data values can be inserted above/below
# them. Inserting 4 blank rows at row 5 moves row5 -> row9 and row6 -> row10.
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).Insert()

# New "Source Type" line (bold + underlined), now at row 7
$ws.Range("A7").Value = "Source Type: Statistical Institution"
$ws.Range("A7").Font.Bold = $true
$ws.Range("A7").Font.Underline = $true

# New data row: percentages for Micro / SMEs / MSMEs, as text, default style
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "30.8"
$ws.Range("B10").Style = "Normal"

$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "64.9"
$ws.Range("C10").Style = "Normal"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "95.7"
$ws.Range("D10").Style = "Normal"

# New source citation line (italic), new row 11
$ws.Range("A11").Value = "Source: SOM - Stat. Office of Montenegro, 2010"
$ws.Range("A11").Font.Italic = $true
